# Updated cryptos list on Fri Jul 19 23:47:04 UTC 2024 with GitHub Actions
# Refresh the Price (D) / Volume(1h) (E) columns with the latest scrape.
# Note: several Price values are plain decimals (e.g. "592.01") that Excel
# would otherwise auto-coerce to numbers (losing the literal-text
# formatting, e.g. trailing zeros). Force those specific cells to Text
# ("@") before assigning so they round-trip exactly as strings, matching
# the source data's inline-string storage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "66.663.80"
$ws.Range("D3").Value = "3.502.22"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.01"
$ws.Range("E5").Value = "  +3.42%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.80"
$ws.Range("E6").Value = "  +5.93%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "3.500.78"
$ws.Range("E8").Value = "  +1.83%  "
$ws.Range("E9").Value = "  +0.05%  "
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.125"
$ws.Range("E11").Value = "  +4.68%  "
$ws.Range("E12").Value = "  +2.98%  "
$ws.Range("D13").Value = "4.112.52"
$ws.Range("E13").Value = "  +2.08%  "
$ws.Range("E14").Value = "  +0.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.10"
$ws.Range("E15").Value = "  +3.32%  "
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("D17").Value = "66.651.38"
$ws.Range("E17").Value = "  +4.02%  "
$ws.Range("D18").Value = "3.535.51"
$ws.Range("E18").Value = "  +2.87%  "
$ws.Range("E19").Value = "  +3.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.03"
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "389.85"
$ws.Range("E21").Value = "  +1.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.97"
$ws.Range("E22").Value = "  +1.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.36"
$ws.Range("E23").Value = "  +2.55%  "
$ws.Range("B24").Value = "Dai"
$ws.Range("C24").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("B25").Value = "PEPE"
$ws.Range("C25").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000127"
$ws.Range("E25").Value = "  +8.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.532"
$ws.Range("E26").Value = "  +3.15%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.11"
$ws.Range("E27").Value = "  +4.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.181"
$ws.Range("E28").Value = "  +1.44%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.42"
$ws.Range("E30").Value = "  +6.24%  "
$ws.Range("E31").Value = "  +4.72%  "
$ws.Range("E32").Value = "  +2.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.57"
$ws.Range("E33").Value = "  +2.25%  "
$ws.Range("E35").Value = "  -0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.62"
$ws.Range("E36").Value = "  +6.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.98"
$ws.Range("E37").Value = "  +0.12%  "
$ws.Range("E38").Value = "  +4.60%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.90"
$ws.Range("E39").Value = "  +3.13%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0746"
$ws.Range("E40").Value = "  +2.60%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "26.64"
$ws.Range("E41").Value = "  +2.06%  "
$ws.Range("E42").Value = "  +5.80%  "
$ws.Range("D43").Value = "2.829.25"
$ws.Range("E43").Value = "  +0.23%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "6.64"
$ws.Range("E44").Value = "  +2.87%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "43.40"
$ws.Range("E45").Value = "  +0.82%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "26.33"
$ws.Range("E47").Value = "  +3.01%  "
$ws.Range("E48").Value = "  +4.08%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "353.87"
$ws.Range("E49").Value = "  +5.58%  "
$ws.Range("E50").Value = "  +2.30%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "34.01"
$ws.Range("E51").Value = "  +14.05%  "
